$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("index")

# Widen download_location column (H) to fit the new longer path text
$ws.Columns.Item(8).ColumnWidth = 84.67

# Row 161 - EQUATES 2002
$ws.Range("A161").Value = "EQUATES"
$ws.Range("B161").Value = 2002
$ws.Range("C161").Value = "smoke_flat_file"
$ws.Range("D161").Value = "onroad"
$ws.Range("I161").Value = "https://drive.usercontent.google.com/download?id=1jw_0216cTDSw-FrojrgGbKU5myl23ybv&export=download&authuser=0&confirm=t&uuid=fccfb764-8849-44fc-b506-7de81bbdd4e8&at=AO7h07e29Tc-CjdVA7IphIHXaT1L:1727277124697"

# Row 162 - EQUATES 2005
$ws.Range("A162").Value = "EQUATES"
$ws.Range("B162").Value = 2005
$ws.Range("C162").Value = "smoke_flat_file"
$ws.Range("D162").Value = "onroad"
$ws.Range("I162").Value = "https://drive.usercontent.google.com/download?id=1W1QBtNB89FfWdSE3RQdnFoedwjgn0ctB&export=download&authuser=0&confirm=t&uuid=9768171e-056f-45b5-aba8-71f4a69146dc&at=AO7h07eZwuW4OJnGLz83u0LU7rHi:1727277818221"

# Row 163 - EQUATES 2008
$ws.Range("A163").Value = "EQUATES"
$ws.Range("B163").Value = 2008
$ws.Range("C163").Value = "smoke_flat_file"
$ws.Range("D163").Value = "onroad"
$ws.Range("I163").Value = "https://drive.usercontent.google.com/download?id=1Y6D0YJuMy97hlhG58R5LjmVrVqucNlgy&export=download&authuser=0&confirm=t&uuid=0bf155e4-d9e9-44a9-ab43-8d93505a6b9b&at=AO7h07eoLdEncM7sZJSEMmhw7bBE:1727277879290"

# Row 165 - EQUATES 2014
$ws.Range("A165").Value = "EQUATES"
$ws.Range("B165").Value = 2014
$ws.Range("C165").Value = "smoke_flat_file"
$ws.Range("D165").Value = "onroad"
$ws.Range("I165").Value = "https://drive.usercontent.google.com/download?id=1TTzPa0vse1_z0On-82uBF53IQ6zY5DIX&export=download&authuser=0&confirm=t&uuid=33810c7d-2a28-4d52-bfc2-68598662f428&at=AO7h07cwKOMiDUnr7tirwMaGCb8t:1727277945324"

# Row 166 - EQUATES 2017
$ws.Range("A166").Value = "EQUATES"
$ws.Range("B166").Value = 2017
$ws.Range("C166").Value = "smoke_flat_file"
$ws.Range("D166").Value = "onroad"
$ws.Range("I166").Value = "https://drive.usercontent.google.com/download?id=1QIaSG3kyReCLUJXXgHubUZKuo1Q5pAsl&export=download&authuser=0&confirm=t&uuid=348fa232-1e5b-445c-9c37-5b49be9819ae&at=AO7h07e7A3qXQiUtgm2L4-zYU3Io:1727277967596"

# Row 167 - EQUATES 2019
$ws.Range("A167").Value = "EQUATES"
$ws.Range("B167").Value = 2019
$ws.Range("C167").Value = "smoke_flat_file"
$ws.Range("D167").Value = "onroad"
$ws.Range("I167").Value = "https://drive.usercontent.google.com/download?id=1ZCAzfSxfzehgKT3RdDO3K6MOko8heOii&export=download&authuser=0&confirm=t&uuid=e597b786-efbb-4922-8161-0af94006d6f0&at=AO7h07cjDkT1iJ9hgSwJ7nIxUY1V:1727277988258"

# Row 164 - EQUATES 2011 (filled last, matching original shared-string order)
$ws.Range("A164").Value = "EQUATES"
$ws.Range("B164").Value = 2011
$ws.Range("C164").Value = "smoke_flat_file"
$ws.Range("D164").Value = "onroad"
$ws.Range("I164").Value = "https://drive.usercontent.google.com/download?id=17qLqwJllqX-XloScSNaxdviQDspTufep&export=download&authuser=0&confirm=t&uuid=d6b82ccf-f5c3-41d2-9b2b-b7bfbf731eeb&at=AO7h07cjoWUwZU0ju9kc4N35ZFQ6:1727278302203"

# download_location column (H) for all 7 new rows, with wrap text
$ws.Range("H161:H167").WrapText = $true
$ws.Range("H161").Value = "_transportation/data-raw/epa/air_emissions_modeling/EQUATES/CMAS_Data_Warehouse"
$ws.Range("H162").Value = "_transportation/data-raw/epa/air_emissions_modeling/EQUATES/CMAS_Data_Warehouse"
$ws.Range("H163").Value = "_transportation/data-raw/epa/air_emissions_modeling/EQUATES/CMAS_Data_Warehouse"
$ws.Range("H164").Value = "_transportation/data-raw/epa/air_emissions_modeling/EQUATES/CMAS_Data_Warehouse"
$ws.Range("H165").Value = "_transportation/data-raw/epa/air_emissions_modeling/EQUATES/CMAS_Data_Warehouse"
$ws.Range("H166").Value = "_transportation/data-raw/epa/air_emissions_modeling/EQUATES/CMAS_Data_Warehouse"
$ws.Range("H167").Value = "_transportation/data-raw/epa/air_emissions_modeling/EQUATES/CMAS_Data_Warehouse"

# Note for 2019 EQUATES row (G column)
$ws.Range("G167").Value = "Contains N2O"

# "base url" sheet - new row 5 with CMAS Data Warehouse entry
$ws2 = $wb.Worksheets.Item("base url")
$ws2.Range("A5").Value = "EQUATES CMAS Data Warehouse"
$ws2.Range("B5").Value = "https://drive.google.com/drive/folders/1G2_LBLy7_n91Ur0ulsLZ9zwGs3luTzn2"
[void]$ws2.Range("B5").Select()

# "notes" sheet - new row 7 note about CMAS data warehouse
$ws3 = $wb.Worksheets.Item("notes")
$ws3.Range("A7").Value = "Data downloaded from the CMAS Data Warehouse Google Drive includes many more sectors and scripts."
[void]$ws3.Range("A8").Select()

# Return focus to "index" sheet and set its selection last so it stays the active tab
[void]$ws.Range("C161:C167").Select()
